# Auto-generated Excel COM-interop script to apply the diff changes
# Updates currentAveragePrice / LevePrice / LeveProfit figures across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 20835522
$ws.Range("J100").Value = 2833.3333
$ws.Range("L100").Value = 2833.3333
$ws.Range("N100").Value = -3915.3333
$ws.Range("H103").Value = 650705
$ws.Range("I103").Value = 612.5
$ws.Range("J103").Value = 1084100
$ws.Range("K103").Value = 1837.5
$ws.Range("L103").Value = 3252300
$ws.Range("M103").Value = -1251.5
$ws.Range("N103").Value = -3253472
$ws.Range("H112").Value = 892.2835700000001
$ws.Range("J112").Value = 934.11865
$ws.Range("L112").Value = 2802.35595
$ws.Range("N112").Value = -5018.35595
$ws.Range("H129").Value = 905.3396
$ws.Range("J129").Value = 939.04083
$ws.Range("L129").Value = 2817.12249
$ws.Range("N129").Value = -12817.12249
$ws.Range("H137").Value = 1998.1923
$ws.Range("I137").Value = 1168.1818
$ws.Range("J137").Value = 2606.8667
$ws.Range("K137").Value = 3504.5454
$ws.Range("L137").Value = 7820.6001
$ws.Range("M137").Value = -954.5454
$ws.Range("N137").Value = -12920.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4182.2173
$ws.Range("I32").Value = 3484.614
$ws.Range("K32").Value = 3484.614
$ws.Range("M32").Value = -3197.614
$ws.Range("H74").Value = 1281
$ws.Range("I74").Value = 1015.6
$ws.Range("K74").Value = 1015.6
$ws.Range("M74").Value = -141.6
$ws.Range("H77").Value = 1281
$ws.Range("I77").Value = 1015.6
$ws.Range("K77").Value = 5078
$ws.Range("M77").Value = -710
$ws.Range("H139").Value = 33877.777
$ws.Range("J139").Value = 33877.777
$ws.Range("L139").Value = 33877.777
$ws.Range("N139").Value = -44157.777

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 59723.816
$ws.Range("J138").Value = 59723.816
$ws.Range("L138").Value = 59723.816
$ws.Range("N138").Value = -70003.81599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2338.026
$ws.Range("I31").Value = 1809.1333
$ws.Range("J31").Value = 2465.984
$ws.Range("K31").Value = 1809.1333
$ws.Range("L31").Value = 2465.984
$ws.Range("M31").Value = -1514.1333
$ws.Range("N31").Value = -3055.984
$ws.Range("H34").Value = 2338.026
$ws.Range("I34").Value = 1809.1333
$ws.Range("J34").Value = 2465.984
$ws.Range("K34").Value = 1809.1333
$ws.Range("L34").Value = 2465.984
$ws.Range("M34").Value = -1607.1333
$ws.Range("N34").Value = -2869.984

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 139066.75
$ws.Range("I5").Value = 12941.25
$ws.Range("J5").Value = 167094.64
$ws.Range("K5").Value = 38823.75
$ws.Range("L5").Value = 501283.92
$ws.Range("M5").Value = -38711.75
$ws.Range("N5").Value = -501507.92
$ws.Range("H97").Value = 10000398
$ws.Range("J97").Value = 493
$ws.Range("L97").Value = 1479
$ws.Range("N97").Value = -2471
$ws.Range("H107").Value = 810.1321
$ws.Range("J107").Value = 1029.5135
$ws.Range("L107").Value = 3088.5405
$ws.Range("N107").Value = -6928.5405
$ws.Range("H113").Value = 154308.78
$ws.Range("I113").Value = 454.6744
$ws.Range("J113").Value = 455023.62
$ws.Range("K113").Value = 1364.0232
$ws.Range("L113").Value = 1365070.86
$ws.Range("M113").Value = 805.9767999999999
$ws.Range("N113").Value = -1369410.86
$ws.Range("H122").Value = 466.94116
$ws.Range("I122").Value = 399.6
$ws.Range("K122").Value = 3596.4
$ws.Range("M122").Value = -1146.4
$ws.Range("H131").Value = 15494154
$ws.Range("J131").Value = 16668001
$ws.Range("L131").Value = 50004003
$ws.Range("N131").Value = -50014083
$ws.Range("H132").Value = 1567136.1
$ws.Range("I132").Value = 651.4286
$ws.Range("J132").Value = 1951886.9
$ws.Range("K132").Value = 5862.8574
$ws.Range("L132").Value = 17566982.1
$ws.Range("M132").Value = -3332.8574
$ws.Range("N132").Value = -17572042.1
$ws.Range("H135").Value = 139066.75
$ws.Range("I135").Value = 12941.25
$ws.Range("J135").Value = 167094.64
$ws.Range("K135").Value = 116471.25
$ws.Range("L135").Value = 1503851.76
$ws.Range("M135").Value = -113936.25
$ws.Range("N135").Value = -1508921.76

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9221.6
$ws.Range("I70").Value = 13504
$ws.Range("J70").Value = 6366.6665
$ws.Range("K70").Value = 13504
$ws.Range("L70").Value = 6366.6665
$ws.Range("M70").Value = -13234
$ws.Range("N70").Value = -6906.6665
$ws.Range("H73").Value = 9221.6
$ws.Range("I73").Value = 13504
$ws.Range("J73").Value = 6366.6665
$ws.Range("K73").Value = 13504
$ws.Range("L73").Value = 6366.6665
$ws.Range("M73").Value = -12568
$ws.Range("N73").Value = -8238.666499999999
$ws.Range("H116").Value = 28333.334
$ws.Range("J116").Value = 28333.334
$ws.Range("L116").Value = 28333.334
$ws.Range("N116").Value = -37511.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2325.3333
$ws.Range("I7").Value = 2000.6154
$ws.Range("J7").Value = 2709.0908
$ws.Range("K7").Value = 2000.6154
$ws.Range("L7").Value = 2709.0908
$ws.Range("M7").Value = -1888.6154
$ws.Range("N7").Value = -2933.0908
$ws.Range("H22").Value = 3369349
$ws.Range("I22").Value = 11112301
$ws.Range("J22").Value = 2847.913
$ws.Range("K22").Value = 11112301
$ws.Range("L22").Value = 2847.913
$ws.Range("M22").Value = -11112006
$ws.Range("N22").Value = -3437.913
$ws.Range("H27").Value = 3369349
$ws.Range("I27").Value = 11112301
$ws.Range("J27").Value = 2847.913
$ws.Range("K27").Value = 11112301
$ws.Range("L27").Value = 2847.913
$ws.Range("M27").Value = -11112194
$ws.Range("N27").Value = -3061.913
$ws.Range("H40").Value = 76927610
$ws.Range("I40").Value = 142860700
$ws.Range("J40").Value = 5667.5
$ws.Range("K40").Value = 142860700
$ws.Range("L40").Value = 5667.5
$ws.Range("M40").Value = -142860564
$ws.Range("N40").Value = -5939.5
$ws.Range("H46").Value = 27779256
$ws.Range("I46").Value = 41667584
$ws.Range("K46").Value = 41667584
$ws.Range("M46").Value = -41667396
$ws.Range("H55").Value = 20833724
$ws.Range("I55").Value = 285.36365
$ws.Range("J55").Value = 38462016
$ws.Range("K55").Value = 285.36365
$ws.Range("L55").Value = 38462016
$ws.Range("M55").Value = -112.36365
$ws.Range("N55").Value = -38462362
$ws.Range("H61").Value = 2632.2222
$ws.Range("J61").Value = 1848.25
$ws.Range("L61").Value = 1848.25
$ws.Range("N61").Value = -2252.25
$ws.Range("H82").Value = 7879352
$ws.Range("I82").Value = 1551.7778
$ws.Range("J82").Value = 18007952
$ws.Range("K82").Value = 1551.7778
$ws.Range("L82").Value = 18007952
$ws.Range("M82").Value = -1190.7778
$ws.Range("N82").Value = -18008674
$ws.Range("H85").Value = 7879352
$ws.Range("I85").Value = 1551.7778
$ws.Range("J85").Value = 18007952
$ws.Range("K85").Value = 1551.7778
$ws.Range("L85").Value = 18007952
$ws.Range("M85").Value = -303.7778000000001
$ws.Range("N85").Value = -18010448
$ws.Range("H100").Value = 2499.75
$ws.Range("I100").Value = 1999.5
$ws.Range("K100").Value = 1999.5
$ws.Range("M100").Value = -1458.5
$ws.Range("H113").Value = 2632.2222
$ws.Range("J113").Value = 1848.25
$ws.Range("L113").Value = 1848.25
$ws.Range("N113").Value = -6188.25
$ws.Range("H126").Value = 2325.3333
$ws.Range("I126").Value = 2000.6154
$ws.Range("J126").Value = 2709.0908
$ws.Range("K126").Value = 6001.8462
$ws.Range("L126").Value = 8127.2724
$ws.Range("M126").Value = -3531.8462
$ws.Range("N126").Value = -13067.2724
$ws.Range("H132").Value = 27087856
$ws.Range("I132").Value = 28893512
$ws.Range("J132").Value = 3005
$ws.Range("K132").Value = 86680536
$ws.Range("L132").Value = 9015
$ws.Range("N132").Value = -14075
$ws.Range("H136").Value = 2828.1191
$ws.Range("I136").Value = 1838.8667
$ws.Range("J136").Value = 5301.25
$ws.Range("K136").Value = 5516.6001
$ws.Range("L136").Value = 15903.75
$ws.Range("M136").Value = -2966.6001
$ws.Range("N136").Value = -21003.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1773.091
$ws.Range("I81").Value = 1614.2858
$ws.Range("J81").Value = 2051
$ws.Range("K81").Value = 3228.5716
$ws.Range("L81").Value = 4102
$ws.Range("M81").Value = -2167.5716
$ws.Range("N81").Value = -6224
$ws.Range("H84").Value = 1773.091
$ws.Range("I84").Value = 1614.2858
$ws.Range("J84").Value = 2051
$ws.Range("K84").Value = 16142.858
$ws.Range("L84").Value = 20510
$ws.Range("M84").Value = -10838.858
$ws.Range("N84").Value = -31118
$ws.Range("H113").Value = 1311.0344
$ws.Range("I113").Value = 767.6
$ws.Range("J113").Value = 1893.2858
$ws.Range("K113").Value = 2302.8
$ws.Range("L113").Value = 5679.857400000001
$ws.Range("M113").Value = -132.8000000000002
$ws.Range("N113").Value = -10019.8574
